$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.455.30"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.570.97"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.27"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3716"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.37"
$ws.Range("E8").Value = "  -3.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3311"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.134"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07487"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.74"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.931"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.889"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.571.39"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001119"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.79"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06739"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.353"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.53"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.05"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.452.83"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.382"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.577"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.79"
$ws.Range("E27").Value = "  +3.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.71"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.028"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.58"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.748.26"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.058"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.122"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.792"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08342"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02462"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2269"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06410"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.290"
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.345"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6329"
$ws.Range("E42").Value = "  +2.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.29"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.87"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6168"
$ws.Range("E45").Value = "  +6.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.772"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.056"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.87"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.211"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.03"
$ws.Range("E51").Value = "  +2.50%  "
